$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4617.087
$ws.Range("I40").Value = 3193.4375
$ws.Range("J40").Value = 7871.143
$ws.Range("K40").Value = 3193.4375
$ws.Range("L40").Value = 7871.143
$ws.Range("M40").Value = -3018.4375
$ws.Range("N40").Value = -8221.143
$ws.Range("H51").Value = 10192.647
$ws.Range("I51").Value = 9562.666999999999
$ws.Range("J51").Value = 10327.643
$ws.Range("K51").Value = 9562.666999999999
$ws.Range("L51").Value = 10327.643
$ws.Range("M51").Value = -9078.666999999999
$ws.Range("N51").Value = -11295.643
$ws.Range("H74").Value = 17823.611
$ws.Range("I74").Value = 19253.75
$ws.Range("J74").Value = 14963.333
$ws.Range("K74").Value = 19253.75
$ws.Range("L74").Value = 14963.333
$ws.Range("M74").Value = -18317.75
$ws.Range("N74").Value = -16835.333
$ws.Range("H77").Value = 17823.611
$ws.Range("I77").Value = 19253.75
$ws.Range("J77").Value = 14963.333
$ws.Range("K77").Value = 96268.75
$ws.Range("L77").Value = 74816.66500000001
$ws.Range("M77").Value = -91588.75
$ws.Range("N77").Value = -84176.66500000001
$ws.Range("H129").Value = 6829
$ws.Range("I129").Value = 7435.6665
$ws.Range("J129").Value = 1369
$ws.Range("K129").Value = 22306.9995
$ws.Range("L129").Value = 4107
$ws.Range("M129").Value = -17306.9995
$ws.Range("N129").Value = -14107
$ws.Range("H131").Value = 1169.5
$ws.Range("I131").Value = 1169.5
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 3508.5
$ws.Range("L131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 3439.9614
$ws.Range("I132").Value = 2884.6667
$ws.Range("K132").Value = 8654.000100000001
$ws.Range("M132").Value = -6124.000100000001
$ws.Range("H137").Value = 3233.7273
$ws.Range("I137").Value = 3197.1
$ws.Range("K137").Value = 9591.299999999999
$ws.Range("M137").Value = -7041.299999999999
$ws.Range("H138").Value = 2367.139
$ws.Range("I138").Value = 1915.2
$ws.Range("J138").Value = 3394.2727
$ws.Range("K138").Value = 5745.6
$ws.Range("L138").Value = 10182.8181
$ws.Range("M138").Value = -605.6000000000004
$ws.Range("N138").Value = -20462.8181

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 955.86365
$ws.Range("I2").Value = 810.2857
$ws.Range("K2").Value = 810.2857
$ws.Range("M2").Value = -697.2857
$ws.Range("H32").Value = 8107.3535
$ws.Range("I32").Value = 4439.8594
$ws.Range("K32").Value = 4439.8594
$ws.Range("M32").Value = -4152.8594
$ws.Range("H74").Value = 3676.2415
$ws.Range("I74").Value = 2232.44
$ws.Range("K74").Value = 2232.44
$ws.Range("M74").Value = -1358.44
$ws.Range("H77").Value = 3676.2415
$ws.Range("I77").Value = 2232.44
$ws.Range("K77").Value = 11162.2
$ws.Range("M77").Value = -6794.200000000001
$ws.Range("H116").Value = 955.86365
$ws.Range("I116").Value = 810.2857
$ws.Range("K116").Value = 810.2857
$ws.Range("M116").Value = 1483.7143
$ws.Range("H139").Value = 101936.75
$ws.Range("J139").Value = 101936.75
$ws.Range("L139").Value = 101936.75
$ws.Range("N139").Value = -112216.75

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 955.86365
$ws.Range("I3").Value = 810.2857
$ws.Range("K3").Value = 810.2857
$ws.Range("M3").Value = -696.2857
$ws.Range("H20").Value = 986.95
$ws.Range("I20").Value = 867.7692
$ws.Range("K20").Value = 867.7692
$ws.Range("M20").Value = -620.7692
$ws.Range("H86").Value = 3301.75
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 3301.75
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -5547.75
$ws.Range("H89").Value = 3301.75
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 3301.75
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -27740.75
$ws.Range("H134").Value = 2675.6562
$ws.Range("I134").Value = 2457.9355
$ws.Range("J134").Value = 9425
$ws.Range("K134").Value = 7373.806500000001
$ws.Range("L134").Value = 28275
$ws.Range("M134").Value = -4838.806500000001
$ws.Range("N134").Value = -33345

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8240.799999999999
$ws.Range("I31").Value = 4281.4287
$ws.Range("J31").Value = 10880.381
$ws.Range("K31").Value = 4281.4287
$ws.Range("L31").Value = 10880.381
$ws.Range("M31").Value = -3986.4287
$ws.Range("N31").Value = -11470.381
$ws.Range("H34").Value = 8240.799999999999
$ws.Range("I34").Value = 4281.4287
$ws.Range("J34").Value = 10880.381
$ws.Range("K34").Value = 4281.4287
$ws.Range("L34").Value = 10880.381
$ws.Range("M34").Value = -4079.4287
$ws.Range("N34").Value = -11284.381
$ws.Range("H58").Value = 1773.6342
$ws.Range("I58").Value = 1150.5
$ws.Range("K58").Value = 1150.5
$ws.Range("M58").Value = -947.5
$ws.Range("H99").Value = 1890.1086
$ws.Range("I99").Value = 1921.25
$ws.Range("J99").Value = 1682.5
$ws.Range("K99").Value = 1921.25
$ws.Range("L99").Value = 1682.5
$ws.Range("M99").Value = -423.25
$ws.Range("N99").Value = -4678.5
$ws.Range("H126").Value = 1890.1086
$ws.Range("I126").Value = 1921.25
$ws.Range("J126").Value = 1682.5
$ws.Range("K126").Value = 5763.75
$ws.Range("L126").Value = 5047.5
$ws.Range("M126").Value = -3293.75
$ws.Range("N126").Value = -9987.5
$ws.Range("H132").Value = 2244.1304
$ws.Range("I132").Value = 1677.7435
$ws.Range("J132").Value = 5399.7144
$ws.Range("K132").Value = 5033.2305
$ws.Range("L132").Value = 16199.1432
$ws.Range("M132").Value = -2503.2305
$ws.Range("N132").Value = -21259.1432
$ws.Range("H136").Value = 1773.6342
$ws.Range("I136").Value = 1150.5
$ws.Range("K136").Value = 3451.5
$ws.Range("M136").Value = -901.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 347.16666
$ws.Range("I44").Value = 310.4
$ws.Range("J44").Value = 531
$ws.Range("K44").Value = 931.1999999999999
$ws.Range("L44").Value = 1593
$ws.Range("M44").Value = -533.1999999999999
$ws.Range("N44").Value = -2389

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 38399.6
$ws.Range("I43").Value = 18000
$ws.Range("J43").Value = 43499.5
$ws.Range("K43").Value = 18000
$ws.Range("L43").Value = 43499.5
$ws.Range("M43").Value = -17849
$ws.Range("N43").Value = -43801.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3404.0557
$ws.Range("I22").Value = 2628.1
$ws.Range("J22").Value = 4374
$ws.Range("K22").Value = 2628.1
$ws.Range("L22").Value = 4374
$ws.Range("M22").Value = -2333.1
$ws.Range("N22").Value = -4964
$ws.Range("H27").Value = 3404.0557
$ws.Range("I27").Value = 2628.1
$ws.Range("J27").Value = 4374
$ws.Range("K27").Value = 2628.1
$ws.Range("L27").Value = 4374
$ws.Range("M27").Value = -2521.1
$ws.Range("N27").Value = -4588

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3251.5217
$ws.Range("I132").Value = 2621.0715
$ws.Range("J132").Value = 4232.222
$ws.Range("K132").Value = 7863.2145
$ws.Range("L132").Value = 12696.666
$ws.Range("M132").Value = -5333.2145
$ws.Range("N132").Value = -17756.666
$ws.Range("H136").Value = 2877.2
$ws.Range("I136").Value = 2707.0715
$ws.Range("J136").Value = 3157.4119
$ws.Range("K136").Value = 8121.2145
$ws.Range("L136").Value = 9472.235700000001
$ws.Range("M136").Value = -5571.2145
$ws.Range("N136").Value = -14572.2357
